# Update share of costs that must be covered to be deemed profitable
# for dispatchable plants on sheet "SoCtMbCtbDP".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SoCtMbCtbDP")

# Update individual plant-type values in column B.
$ws.Range("B2").Value = 1      # hard coal
$ws.Range("B3").Value = 1      # natural gas steam turbine
$ws.Range("B5").Value = 1      # nuclear
$ws.Range("B14").Value = 0.9   # lignite
$ws.Range("B18").Value = 1     # municipal solid waste
$ws.Range("B19").Value = 1     # hard coal w CCS
$ws.Range("B20").Value = 1     # natural gas combined cycle w CCS
$ws.Range("B21").Value = 1     # biomass w CCS
$ws.Range("B22").Value = 1     # lignite w CCS
$ws.Range("B23").Value = 1     # small modular reactor
$ws.Range("B24").Value = 1     # hydrogen combustion turbine
$ws.Range("B25").Value = 1     # hydrogen combined cycle

# Reflect the final selection state captured in the saved workbook.
$ws.Activate()
$ws.Range("B18:B25").Select()
